$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPG3JQWOZ%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC0aCXVzLXdlc3QtMiJIMEYCIQC1zBvgPHY5FgK2J7oy01XgLuo1dPTM2TAbec907QphQgIhANw6mwD8rdVkc6JumldHQyDzZ8G4vZ5TrHE%2FtZFAcik9KuYDCJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyWY4ZJoAs2Qw3n9dwqugM2nceXvnLe8Ai%2ByEMTNYs1yKVi5z9KyN6j5Z%2BDOA5UfkSR5evrLDL00bgNXYft7xqNheZE6mklbP3guRSzNzJORQJ7zfRN3XwXhpXy%2FiPLc%2BcNw1BdTBEZCrj4OfD4VCMQ07Ob%2F70ybNkjY8g8X3CfmNvcY3iM2nBjYdSnxnsqqozC1FxXQtEsseJ7dB7jWoSkMafnhK4Z3XivTC%2FuvDd3ahenJBXvDk9%2Fq3N95oC5y6pKjZ9qd9ZNvfl4Xlub6c%2BNlEUFr0ulFvvXuxSoVeqRNqz1PIDH%2BGhYNaTEvWIazj%2B%2FcV%2Bkryz81Z47MstQk1Q4YadrObzQNHy8fMJgrRyD%2BiA4Qpfj41dEas3yh%2Ft8B6GUXEBE1evoXqtLbAsqFupwjGEWKZnOueKxurWYfEI%2F6bCLw%2FDkz9KqBHunqrYpJ3l2IxGhmN0kFrjk%2BUklGtNUCi%2FDGDQdHa70uRo2sY0Xt%2F4zEHm2bgY1UjbCGfZDPCUKKpPbedYBcFjU%2FR9RCb5MZdo%2FodI%2B8Hrq9SVXAkpCHpWfrwjcNblBVjugEluXQ5S%2BPwt5GZ5urvNfJKsXBMjPZ8387Ham4tftMMrlpr8GOqQB4J6guDqUtPRubaoSfZXItv2OJEz%2F6xov%2B0hIFx%2BipHw9V3jWvd41zGFo4FNVYyX07g4CXtCpO3ci%2FZ2qTPUA9FpnB52%2Bn3UxhyH7prQeBbfbaUdf6PFrkrlRSKQXWpZ0ND7xnULw2LlBBn1dGIH9gy5HkY7ZV7CVLHKiQ6ERmf52t3WSdF%2FGuWne%2B2PZS7Fin7pKtjh%2FJoCrp05IXBhFOas9ATc%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=d5b5c2281291299747437448761bfcb9edebdf4193487f61c7ad22ed5005025f'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPG3JQWOZ%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC0aCXVzLXdlc3QtMiJIMEYCIQC1zBvgPHY5FgK2J7oy01XgLuo1dPTM2TAbec907QphQgIhANw6mwD8rdVkc6JumldHQyDzZ8G4vZ5TrHE%2FtZFAcik9KuYDCJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyWY4ZJoAs2Qw3n9dwqugM2nceXvnLe8Ai%2ByEMTNYs1yKVi5z9KyN6j5Z%2BDOA5UfkSR5evrLDL00bgNXYft7xqNheZE6mklbP3guRSzNzJORQJ7zfRN3XwXhpXy%2FiPLc%2BcNw1BdTBEZCrj4OfD4VCMQ07Ob%2F70ybNkjY8g8X3CfmNvcY3iM2nBjYdSnxnsqqozC1FxXQtEsseJ7dB7jWoSkMafnhK4Z3XivTC%2FuvDd3ahenJBXvDk9%2Fq3N95oC5y6pKjZ9qd9ZNvfl4Xlub6c%2BNlEUFr0ulFvvXuxSoVeqRNqz1PIDH%2BGhYNaTEvWIazj%2B%2FcV%2Bkryz81Z47MstQk1Q4YadrObzQNHy8fMJgrRyD%2BiA4Qpfj41dEas3yh%2Ft8B6GUXEBE1evoXqtLbAsqFupwjGEWKZnOueKxurWYfEI%2F6bCLw%2FDkz9KqBHunqrYpJ3l2IxGhmN0kFrjk%2BUklGtNUCi%2FDGDQdHa70uRo2sY0Xt%2F4zEHm2bgY1UjbCGfZDPCUKKpPbedYBcFjU%2FR9RCb5MZdo%2FodI%2B8Hrq9SVXAkpCHpWfrwjcNblBVjugEluXQ5S%2BPwt5GZ5urvNfJKsXBMjPZ8387Ham4tftMMrlpr8GOqQB4J6guDqUtPRubaoSfZXItv2OJEz%2F6xov%2B0hIFx%2BipHw9V3jWvd41zGFo4FNVYyX07g4CXtCpO3ci%2FZ2qTPUA9FpnB52%2Bn3UxhyH7prQeBbfbaUdf6PFrkrlRSKQXWpZ0ND7xnULw2LlBBn1dGIH9gy5HkY7ZV7CVLHKiQ6ERmf52t3WSdF%2FGuWne%2B2PZS7Fin7pKtjh%2FJoCrp05IXBhFOas9ATc%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=ac313bd15ff1eb12177c4416346826c4ca29d87b415b73f5d34f51878877d46c'
$ws.Range('K4').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPG3JQWOZ%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC0aCXVzLXdlc3QtMiJIMEYCIQC1zBvgPHY5FgK2J7oy01XgLuo1dPTM2TAbec907QphQgIhANw6mwD8rdVkc6JumldHQyDzZ8G4vZ5TrHE%2FtZFAcik9KuYDCJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyWY4ZJoAs2Qw3n9dwqugM2nceXvnLe8Ai%2ByEMTNYs1yKVi5z9KyN6j5Z%2BDOA5UfkSR5evrLDL00bgNXYft7xqNheZE6mklbP3guRSzNzJORQJ7zfRN3XwXhpXy%2FiPLc%2BcNw1BdTBEZCrj4OfD4VCMQ07Ob%2F70ybNkjY8g8X3CfmNvcY3iM2nBjYdSnxnsqqozC1FxXQtEsseJ7dB7jWoSkMafnhK4Z3XivTC%2FuvDd3ahenJBXvDk9%2Fq3N95oC5y6pKjZ9qd9ZNvfl4Xlub6c%2BNlEUFr0ulFvvXuxSoVeqRNqz1PIDH%2BGhYNaTEvWIazj%2B%2FcV%2Bkryz81Z47MstQk1Q4YadrObzQNHy8fMJgrRyD%2BiA4Qpfj41dEas3yh%2Ft8B6GUXEBE1evoXqtLbAsqFupwjGEWKZnOueKxurWYfEI%2F6bCLw%2FDkz9KqBHunqrYpJ3l2IxGhmN0kFrjk%2BUklGtNUCi%2FDGDQdHa70uRo2sY0Xt%2F4zEHm2bgY1UjbCGfZDPCUKKpPbedYBcFjU%2FR9RCb5MZdo%2FodI%2B8Hrq9SVXAkpCHpWfrwjcNblBVjugEluXQ5S%2BPwt5GZ5urvNfJKsXBMjPZ8387Ham4tftMMrlpr8GOqQB4J6guDqUtPRubaoSfZXItv2OJEz%2F6xov%2B0hIFx%2BipHw9V3jWvd41zGFo4FNVYyX07g4CXtCpO3ci%2FZ2qTPUA9FpnB52%2Bn3UxhyH7prQeBbfbaUdf6PFrkrlRSKQXWpZ0ND7xnULw2LlBBn1dGIH9gy5HkY7ZV7CVLHKiQ6ERmf52t3WSdF%2FGuWne%2B2PZS7Fin7pKtjh%2FJoCrp05IXBhFOas9ATc%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=bfb04a9b0af4333405fda669625fd008dd042acac30939661084fd4a931ac4a2'
$ws.Range('K5').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPG3JQWOZ%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC0aCXVzLXdlc3QtMiJIMEYCIQC1zBvgPHY5FgK2J7oy01XgLuo1dPTM2TAbec907QphQgIhANw6mwD8rdVkc6JumldHQyDzZ8G4vZ5TrHE%2FtZFAcik9KuYDCJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyWY4ZJoAs2Qw3n9dwqugM2nceXvnLe8Ai%2ByEMTNYs1yKVi5z9KyN6j5Z%2BDOA5UfkSR5evrLDL00bgNXYft7xqNheZE6mklbP3guRSzNzJORQJ7zfRN3XwXhpXy%2FiPLc%2BcNw1BdTBEZCrj4OfD4VCMQ07Ob%2F70ybNkjY8g8X3CfmNvcY3iM2nBjYdSnxnsqqozC1FxXQtEsseJ7dB7jWoSkMafnhK4Z3XivTC%2FuvDd3ahenJBXvDk9%2Fq3N95oC5y6pKjZ9qd9ZNvfl4Xlub6c%2BNlEUFr0ulFvvXuxSoVeqRNqz1PIDH%2BGhYNaTEvWIazj%2B%2FcV%2Bkryz81Z47MstQk1Q4YadrObzQNHy8fMJgrRyD%2BiA4Qpfj41dEas3yh%2Ft8B6GUXEBE1evoXqtLbAsqFupwjGEWKZnOueKxurWYfEI%2F6bCLw%2FDkz9KqBHunqrYpJ3l2IxGhmN0kFrjk%2BUklGtNUCi%2FDGDQdHa70uRo2sY0Xt%2F4zEHm2bgY1UjbCGfZDPCUKKpPbedYBcFjU%2FR9RCb5MZdo%2FodI%2B8Hrq9SVXAkpCHpWfrwjcNblBVjugEluXQ5S%2BPwt5GZ5urvNfJKsXBMjPZ8387Ham4tftMMrlpr8GOqQB4J6guDqUtPRubaoSfZXItv2OJEz%2F6xov%2B0hIFx%2BipHw9V3jWvd41zGFo4FNVYyX07g4CXtCpO3ci%2FZ2qTPUA9FpnB52%2Bn3UxhyH7prQeBbfbaUdf6PFrkrlRSKQXWpZ0ND7xnULw2LlBBn1dGIH9gy5HkY7ZV7CVLHKiQ6ERmf52t3WSdF%2FGuWne%2B2PZS7Fin7pKtjh%2FJoCrp05IXBhFOas9ATc%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=ac27554be1dc2ff3c38498c064b6f641f8a8d65b83b1dc74a61ba1e1d755a977'
$ws.Range('K6').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPG3JQWOZ%2F20250331%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250331T030004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEC0aCXVzLXdlc3QtMiJIMEYCIQC1zBvgPHY5FgK2J7oy01XgLuo1dPTM2TAbec907QphQgIhANw6mwD8rdVkc6JumldHQyDzZ8G4vZ5TrHE%2FtZFAcik9KuYDCJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgyWY4ZJoAs2Qw3n9dwqugM2nceXvnLe8Ai%2ByEMTNYs1yKVi5z9KyN6j5Z%2BDOA5UfkSR5evrLDL00bgNXYft7xqNheZE6mklbP3guRSzNzJORQJ7zfRN3XwXhpXy%2FiPLc%2BcNw1BdTBEZCrj4OfD4VCMQ07Ob%2F70ybNkjY8g8X3CfmNvcY3iM2nBjYdSnxnsqqozC1FxXQtEsseJ7dB7jWoSkMafnhK4Z3XivTC%2FuvDd3ahenJBXvDk9%2Fq3N95oC5y6pKjZ9qd9ZNvfl4Xlub6c%2BNlEUFr0ulFvvXuxSoVeqRNqz1PIDH%2BGhYNaTEvWIazj%2B%2FcV%2Bkryz81Z47MstQk1Q4YadrObzQNHy8fMJgrRyD%2BiA4Qpfj41dEas3yh%2Ft8B6GUXEBE1evoXqtLbAsqFupwjGEWKZnOueKxurWYfEI%2F6bCLw%2FDkz9KqBHunqrYpJ3l2IxGhmN0kFrjk%2BUklGtNUCi%2FDGDQdHa70uRo2sY0Xt%2F4zEHm2bgY1UjbCGfZDPCUKKpPbedYBcFjU%2FR9RCb5MZdo%2FodI%2B8Hrq9SVXAkpCHpWfrwjcNblBVjugEluXQ5S%2BPwt5GZ5urvNfJKsXBMjPZ8387Ham4tftMMrlpr8GOqQB4J6guDqUtPRubaoSfZXItv2OJEz%2F6xov%2B0hIFx%2BipHw9V3jWvd41zGFo4FNVYyX07g4CXtCpO3ci%2FZ2qTPUA9FpnB52%2Bn3UxhyH7prQeBbfbaUdf6PFrkrlRSKQXWpZ0ND7xnULw2LlBBn1dGIH9gy5HkY7ZV7CVLHKiQ6ERmf52t3WSdF%2FGuWne%2B2PZS7Fin7pKtjh%2FJoCrp05IXBhFOas9ATc%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2011%3A00%3A04%20GMT&X-Amz-Signature=18e73a968f94fdafa669cf15d276e7ac69fcb24f5ba2a19e3e08a6505ac8b714'
